$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.391867
$ws.Range("H2").Value = 40.175601
$ws.Range("I2").Value = 0.1577826877470924
$ws.Range("J2").Value = 0.1577826877470924
$ws.Range("M2").Value = 20.29529466666667
$ws.Range("N2").Value = 60.885884
$ws.Range("O2").Value = 0.4032332285476398
$ws.Range("P2").Value = 0.4032332285476398
$ws.Range("Q2").Value = 271.7918869018093
$ws.Range("R2").Value = 2446.126982116284
$ws.Range("S2").Value = 0.06362322258918419
$ws.Range("T2").Value = 0.06362322258918418
$ws.Range("G3").Value = 13.391867
$ws.Range("H3").Value = 40.175601
$ws.Range("I3").Value = 0.1577826877470924
$ws.Range("J3").Value = 0.1577826877470924
$ws.Range("O3").Value = 0.1953894087318433
$ws.Range("P3").Value = 0.1953894087318433
$ws.Range("Q3").Value = 131.6986109282967
$ws.Range("R3").Value = 1185.28749835467
$ws.Range("S3").Value = 0.03082906606702544
$ws.Range("T3").Value = 0.03082906606702543
$ws.Range("G4").Value = 13.391867
$ws.Range("H4").Value = 40.175601
$ws.Range("I4").Value = 0.1577826877470924
$ws.Range("J4").Value = 0.1577826877470924
$ws.Range("M4").Value = 11.81535133333333
$ws.Range("N4").Value = 35.446054
$ws.Range("O4").Value = 0.2347510761885954
$ws.Range("P4").Value = 0.2347510761885954
$ws.Range("Q4").Value = 158.2296136142726
$ws.Range("R4").Value = 1424.066522528454
$ws.Range("S4").Value = 0.03703965575255905
$ws.Range("T4").Value = 0.03703965575255904
$ws.Range("G5").Value = 13.391867
$ws.Range("H5").Value = 40.175601
$ws.Range("I5").Value = 0.1577826877470924
$ws.Range("J5").Value = 0.1577826877470924
$ws.Range("M5").Value = 8.386535
$ws.Range("N5").Value = 25.159605
$ws.Range("O5").Value = 0.1666262865319216
$ws.Range("P5").Value = 0.1666262865319216
$ws.Range("Q5").Value = 112.311361310845
$ws.Range("R5").Value = 1010.802251797605
$ws.Range("S5").Value = 0.02629074333832374
$ws.Range("T5").Value = 0.02629074333832373
$ws.Range("I6").Value = 0.1864313654770604
$ws.Range("J6").Value = 0.1864313654770604
$ws.Range("M6").Value = 20.29529466666667
$ws.Range("N6").Value = 60.885884
$ws.Range("O6").Value = 0.4032332285476398
$ws.Range("P6").Value = 0.4032332285476398
$ws.Range("Q6").Value = 321.1412691987484
$ws.Range("R6").Value = 2890.271422788736
$ws.Range("S6").Value = 0.07517532140386003
$ws.Range("T6").Value = 0.07517532140386003
$ws.Range("I7").Value = 0.1864313654770604
$ws.Range("J7").Value = 0.1864313654770604
$ws.Range("O7").Value = 0.1953894087318433
$ws.Range("P7").Value = 0.1953894087318433
$ws.Range("S7").Value = 0.036426714269633
$ws.Range("T7").Value = 0.036426714269633
$ws.Range("I8").Value = 0.1864313654770604
$ws.Range("J8").Value = 0.1864313654770604
$ws.Range("M8").Value = 11.81535133333333
$ws.Range("N8").Value = 35.446054
$ws.Range("O8").Value = 0.2347510761885954
$ws.Range("P8").Value = 0.2347510761885954
$ws.Range("Q8").Value = 186.9594398867129
$ws.Range("R8").Value = 1682.634958980416
$ws.Range("S8").Value = 0.04376496368104926
$ws.Range("T8").Value = 0.04376496368104926
$ws.Range("I9").Value = 0.1864313654770604
$ws.Range("J9").Value = 0.1864313654770604
$ws.Range("M9").Value = 8.386535
$ws.Range("N9").Value = 25.159605
$ws.Range("O9").Value = 0.1666262865319216
$ws.Range("P9").Value = 0.1666262865319216
$ws.Range("Q9").Value = 132.7037886522133
$ws.Range("R9").Value = 1194.33409786992
$ws.Range("S9").Value = 0.03106436612251805
$ws.Range("T9").Value = 0.03106436612251805
$ws.Range("G10").Value = 7.803333333333334
$ws.Range("H10").Value = 23.41
$ws.Range("I10").Value = 0.09193870479148361
$ws.Range("J10").Value = 0.0919387047914836
$ws.Range("M10").Value = 20.29529466666667
$ws.Range("N10").Value = 60.885884
$ws.Range("O10").Value = 0.4032332285476398
$ws.Range("P10").Value = 0.4032332285476398
$ws.Range("Q10").Value = 158.3709493822222
$ws.Range("R10").Value = 1425.33854444
$ws.Range("S10").Value = 0.0370727407615583
$ws.Range("T10").Value = 0.03707274076155829
$ws.Range("G11").Value = 7.803333333333334
$ws.Range("H11").Value = 23.41
$ws.Range("I11").Value = 0.09193870479148361
$ws.Range("J11").Value = 0.0919387047914836
$ws.Range("O11").Value = 0.1953894087318433
$ws.Range("P11").Value = 0.1953894087318433
$ws.Range("Q11").Value = 76.73972274444445
$ws.Range("R11").Value = 690.6575047
$ws.Range("S11").Value = 0.01796384916877947
$ws.Range("T11").Value = 0.01796384916877947
$ws.Range("G12").Value = 7.803333333333334
$ws.Range("H12").Value = 23.41
$ws.Range("I12").Value = 0.09193870479148361
$ws.Range("J12").Value = 0.0919387047914836
$ws.Range("M12").Value = 11.81535133333333
$ws.Range("N12").Value = 35.446054
$ws.Range("O12").Value = 0.2347510761885954
$ws.Range("P12").Value = 0.2347510761885954
$ws.Range("Q12").Value = 92.19912490444445
$ws.Range("R12").Value = 829.7921241399999
$ws.Range("S12").Value = 0.02158270989318635
$ws.Range("T12").Value = 0.02158270989318634
$ws.Range("G13").Value = 7.803333333333334
$ws.Range("H13").Value = 23.41
$ws.Range("I13").Value = 0.09193870479148361
$ws.Range("J13").Value = 0.0919387047914836
$ws.Range("M13").Value = 8.386535
$ws.Range("N13").Value = 25.159605
$ws.Range("O13").Value = 0.1666262865319216
$ws.Range("P13").Value = 0.1666262865319216
$ws.Range("Q13").Value = 65.44292811666668
$ws.Range("R13").Value = 588.9863530499999
$ws.Range("S13").Value = 0.0153194049679595
$ws.Range("T13").Value = 0.0153194049679595
$ws.Range("G14").Value = 47.856754
$ws.Range("H14").Value = 143.570262
$ws.Range("I14").Value = 0.5638472419843638
$ws.Range("J14").Value = 0.5638472419843636
$ws.Range("M14").Value = 20.29529466666667
$ws.Range("N14").Value = 60.885884
$ws.Range("O14").Value = 0.4032332285476398
$ws.Range("P14").Value = 0.4032332285476398
$ws.Range("Q14").Value = 971.2669242201786
$ws.Range("R14").Value = 8741.402317981609
$ws.Range("S14").Value = 0.2273619437930373
$ws.Range("T14").Value = 0.2273619437930373
$ws.Range("G15").Value = 47.856754
$ws.Range("H15").Value = 143.570262
$ws.Range("I15").Value = 0.5638472419843638
$ws.Range("J15").Value = 0.5638472419843636
$ws.Range("O15").Value = 0.1953894087318433
$ws.Range("P15").Value = 0.1953894087318433
$ws.Range("Q15").Value = 470.6340068443933
$ws.Range("R15").Value = 4235.70606159954
$ws.Range("S15").Value = 0.1101697792264054
$ws.Range("T15").Value = 0.1101697792264054
$ws.Range("G16").Value = 47.856754
$ws.Range("H16").Value = 143.570262
$ws.Range("I16").Value = 0.5638472419843638
$ws.Range("J16").Value = 0.5638472419843636
$ws.Range("M16").Value = 11.81535133333333
$ws.Range("N16").Value = 35.446054
$ws.Range("O16").Value = 0.2347510761885954
$ws.Range("P16").Value = 0.2347510761885954
$ws.Range("Q16").Value = 565.4443621829053
$ws.Range("R16").Value = 5088.999259646147
$ws.Range("S16").Value = 0.1323637468618007
$ws.Range("T16").Value = 0.1323637468618007
$ws.Range("G17").Value = 47.856754
$ws.Range("H17").Value = 143.570262
$ws.Range("I17").Value = 0.5638472419843638
$ws.Range("J17").Value = 0.5638472419843636
$ws.Range("M17").Value = 8.386535
$ws.Range("N17").Value = 25.159605
$ws.Range("O17").Value = 0.1666262865319216
$ws.Range("P17").Value = 0.1666262865319216
$ws.Range("Q17").Value = 401.35234240739
$ws.Range("R17").Value = 3612.17108166651
$ws.Range("S17").Value = 0.09395177210312033
$ws.Range("T17").Value = 0.09395177210312029
